$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 256.875
$ws.Range("I12").Value = 132.6
$ws.Range("J12").Value = 464
$ws.Range("K12").Value = 132.6
$ws.Range("L12").Value = 464
$ws.Range("M12").Value = 37.40000000000001
$ws.Range("N12").Value = -804

$ws.Range("H37").Value = 3000
$ws.Range("J37").Value = 3000
$ws.Range("L37").Value = 9000
$ws.Range("N37").Value = -9252

$ws.Range("H51").Value = 2360.1
$ws.Range("I51").Value = 2240.2
$ws.Range("J51").Value = 2480
$ws.Range("K51").Value = 2240.2
$ws.Range("L51").Value = 2480
$ws.Range("M51").Value = -1756.2
$ws.Range("N51").Value = -3448

$ws.Range("H58").Value = 2592
$ws.Range("J58").Value = 3322.8
$ws.Range("L58").Value = 9968.400000000001
$ws.Range("N58").Value = -10268.4

$ws.Range("H63").Value = 38271
$ws.Range("J63").Value = 38271
$ws.Range("L63").Value = 38271
$ws.Range("N63").Value = -39519

$ws.Range("H66").Value = 38271
$ws.Range("J66").Value = 38271
$ws.Range("L66").Value = 114813
$ws.Range("N66").Value = -121053

$ws.Range("H112").Value = 5564.0586
$ws.Range("J112").Value = 1530.5625
$ws.Range("L112").Value = 4591.6875
$ws.Range("N112").Value = -6807.6875

$ws.Range("H129").Value = 961.8409
$ws.Range("J129").Value = 983.9761999999999
$ws.Range("L129").Value = 2951.9286
$ws.Range("N129").Value = -12951.9286

$ws.Range("H138").Value = 2485.5305
$ws.Range("I138").Value = 1224.6957
$ws.Range("J138").Value = 3600.8845
$ws.Range("K138").Value = 3674.0871
$ws.Range("L138").Value = 10802.6535
$ws.Range("M138").Value = 1465.9129
$ws.Range("N138").Value = -21082.6535

$ws.Range("H139").Value = 65000
$ws.Range("J139").Value = 65000
$ws.Range("L139").Value = 65000
$ws.Range("N139").Value = -75280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9125.223
$ws.Range("I61").Value = 9001
$ws.Range("J61").Value = 9373.666999999999
$ws.Range("K61").Value = 9001
$ws.Range("L61").Value = 9373.666999999999
$ws.Range("M61").Value = -8789
$ws.Range("N61").Value = -9797.666999999999

$ws.Range("H63").Value = 2200.818
$ws.Range("I63").Value = 2200.818
$ws.Range("K63").Value = 2200.818
$ws.Range("M63").Value = -1514.818

$ws.Range("H66").Value = 2200.818
$ws.Range("I66").Value = 2200.818
$ws.Range("K66").Value = 11004.09
$ws.Range("M66").Value = -7572.09

$ws.Range("H136").Value = 9125.223
$ws.Range("I136").Value = 9001
$ws.Range("J136").Value = 9373.666999999999
$ws.Range("K136").Value = 27003
$ws.Range("L136").Value = 28121.001
$ws.Range("M136").Value = -24453
$ws.Range("N136").Value = -33221.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 641796.1
$ws.Range("I31").Value = 5303.3477
$ws.Range("J31").Value = 1060062.8
$ws.Range("K31").Value = 5303.3477
$ws.Range("L31").Value = 1060062.8
$ws.Range("M31").Value = -5008.3477
$ws.Range("N31").Value = -1060652.8

$ws.Range("H34").Value = 641796.1
$ws.Range("I34").Value = 5303.3477
$ws.Range("J34").Value = 1060062.8
$ws.Range("K34").Value = 5303.3477
$ws.Range("L34").Value = 1060062.8
$ws.Range("M34").Value = -5101.3477
$ws.Range("N34").Value = -1060466.8

$ws.Range("H59").Value = 23348.732
$ws.Range("J59").Value = 24202.076
$ws.Range("L59").Value = 24202.076
$ws.Range("N59").Value = -26492.076

$ws.Range("H132").Value = 4815.6113
$ws.Range("I132").Value = 4821.769
$ws.Range("K132").Value = 14465.307
$ws.Range("M132").Value = -11935.307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7251030
$ws.Range("I5").Value = 345.875
$ws.Range("J5").Value = 23824022
$ws.Range("K5").Value = 1037.625
$ws.Range("L5").Value = 71472066
$ws.Range("M5").Value = -925.625
$ws.Range("N5").Value = -71472290

$ws.Range("H22").Value = 125002730
$ws.Range("I22").Value = 333334600
$ws.Range("K22").Value = 1000003800
$ws.Range("M22").Value = -1000003631

$ws.Range("H27").Value = 125002730
$ws.Range("I27").Value = 333334600
$ws.Range("K27").Value = 1000003800
$ws.Range("M27").Value = -1000003698

$ws.Range("H38").Value = 64
$ws.Range("I38").Value = 26.461538
$ws.Range("J38").Value = 104.666664
$ws.Range("K38").Value = 79.384614
$ws.Range("L38").Value = 313.999992
$ws.Range("M38").Value = 267.615386
$ws.Range("N38").Value = -1007.999992

$ws.Range("H68").Value = 2225.8027
$ws.Range("I68").Value = 830.4737
$ws.Range("J68").Value = 3621.1316
$ws.Range("K68").Value = 2491.4211
$ws.Range("L68").Value = 10863.3948
$ws.Range("M68").Value = -1680.4211
$ws.Range("N68").Value = -12485.3948

$ws.Range("H71").Value = 2225.8027
$ws.Range("I71").Value = 830.4737
$ws.Range("J71").Value = 3621.1316
$ws.Range("K71").Value = 7474.263300000001
$ws.Range("L71").Value = 32590.1844
$ws.Range("M71").Value = -3418.263300000001
$ws.Range("N71").Value = -40702.1844

$ws.Range("H80").Value = 2200
$ws.Range("I80").Value = 1500
$ws.Range("K80").Value = 4500
$ws.Range("M80").Value = -3564

$ws.Range("H83").Value = 2200
$ws.Range("I83").Value = 1500
$ws.Range("K83").Value = 13500
$ws.Range("M83").Value = -8820

$ws.Range("H107").Value = 620.1568600000001
$ws.Range("I107").Value = 282.45456
$ws.Range("J107").Value = 2742.8572
$ws.Range("K107").Value = 847.36368
$ws.Range("L107").Value = 8228.571599999999
$ws.Range("M107").Value = 1072.63632
$ws.Range("N107").Value = -12068.5716

$ws.Range("H113").Value = 407.79166
$ws.Range("I113").Value = 473.69232
$ws.Range("J113").Value = 362.70175
$ws.Range("K113").Value = 1421.07696
$ws.Range("L113").Value = 1088.10525
$ws.Range("M113").Value = 748.9230400000001
$ws.Range("N113").Value = -5428.105250000001

$ws.Range("H122").Value = 698.0909
$ws.Range("I122").Value = 243.875
$ws.Range("J122").Value = 957.6429000000001
$ws.Range("K122").Value = 2194.875
$ws.Range("L122").Value = 8618.786100000001
$ws.Range("M122").Value = 255.125
$ws.Range("N122").Value = -13518.7861

$ws.Range("H135").Value = 7251030
$ws.Range("I135").Value = 345.875
$ws.Range("J135").Value = 23824022
$ws.Range("K135").Value = 3112.875
$ws.Range("L135").Value = 214416198
$ws.Range("M135").Value = -577.875
$ws.Range("N135").Value = -214421268

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 22400
$ws.Range("I46").Value = 28500
$ws.Range("J46").Value = 21180
$ws.Range("K46").Value = 28500
$ws.Range("L46").Value = 21180
$ws.Range("M46").Value = -28344
$ws.Range("N46").Value = -21492

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents() | Out-Null

$ws.Range("H107").Value = 609.8
$ws.Range("I107").Value = 184.33333
$ws.Range("J107").Value = 1248
$ws.Range("K107").Value = 184.33333
$ws.Range("L107").Value = 1248
$ws.Range("M107").Value = 1735.66667
$ws.Range("N107").Value = -5088

$ws.Range("H132").Value = 3196.5334
$ws.Range("I132").Value = 2706
$ws.Range("J132").Value = 3625.75
$ws.Range("K132").Value = 8118
$ws.Range("L132").Value = 10877.25
$ws.Range("M132").Value = -5588
$ws.Range("N132").Value = -15937.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 668.7273
$ws.Range("I93").Value = 668.7273
$ws.Range("K93").Value = 668.7273
$ws.Range("M93").Value = 579.2727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1512.8334
$ws.Range("I96").Value = 1490
$ws.Range("J96").Value = 1524.25
$ws.Range("K96").Value = 1490
$ws.Range("L96").Value = 1524.25
$ws.Range("M96").Value = -117
$ws.Range("N96").Value = -4270.25

$ws.Range("H132").Value = 3168.125
$ws.Range("I132").Value = 2882.611
$ws.Range("K132").Value = 8647.832999999999
$ws.Range("M132").Value = -6117.832999999999
